$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 199.64285
$ws.Range("I5").Value = 208.07692
$ws.Range("K5").Value = 208.07692
$ws.Range("M5").Value = -93.07692
$ws.Range("H9").Value = 121.25
$ws.Range("I9").Value = 95
$ws.Range("K9").Value = 95
$ws.Range("M9").Value = 74
$ws.Range("H15").Value = 213803.86
$ws.Range("I15").Value = 213803.86
$ws.Range("K15").Value = 641411.58
$ws.Range("M15").Value = -641242.58
$ws.Range("H92").Value = 267.29413
$ws.Range("I92").Value = 236.26666
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 236.26666
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = 1011.73334
$ws.Range("N92").Value = -2996
$ws.Range("H100").Value = 14307333
$ws.Range("J100").Value = 1252601.6
$ws.Range("L100").Value = 1252601.6
$ws.Range("N100").Value = -1253683.6
$ws.Range("H137").Value = 1708.3549
$ws.Range("I137").Value = 1713.9231
$ws.Range("K137").Value = 5141.7693
$ws.Range("M137").Value = -2591.7693
$ws.Range("H138").Value = 12051398
$ws.Range("I138").Value = 1782.9231
$ws.Range("J138").Value = 17547714
$ws.Range("K138").Value = 5348.7693
$ws.Range("L138").Value = 52643142
$ws.Range("M138").Value = -208.7692999999999
$ws.Range("N138").Value = -52653422

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3172.8901
$ws.Range("I32").Value = 2021.439
$ws.Range("K32").Value = 2021.439
$ws.Range("M32").Value = -1734.439
$ws.Range("H61").Value = 4010.2712
$ws.Range("I61").Value = 5068.4414
$ws.Range("K61").Value = 5068.4414
$ws.Range("M61").Value = -4856.4414
$ws.Range("H74").Value = 2285.5789
$ws.Range("I74").Value = 2200.7058
$ws.Range("K74").Value = 2200.7058
$ws.Range("M74").Value = -1326.7058
$ws.Range("H77").Value = 2285.5789
$ws.Range("I77").Value = 2200.7058
$ws.Range("K77").Value = 11003.529
$ws.Range("M77").Value = -6635.529
$ws.Range("H97").Value = 1212
$ws.Range("I97").Value = 997.1818
$ws.Range("J97").Value = 1999.6666
$ws.Range("K97").Value = 997.1818
$ws.Range("L97").Value = 1999.6666
$ws.Range("M97").Value = -501.1818
$ws.Range("N97").Value = -2991.6666
$ws.Range("H102").Value = 2361.111
$ws.Range("I102").Value = 2318.5715
$ws.Range("J102").Value = 2510
$ws.Range("K102").Value = 2318.5715
$ws.Range("L102").Value = 2510
$ws.Range("M102").Value = -696.5715
$ws.Range("N102").Value = -5754
$ws.Range("H110").Value = 76219.75
$ws.Range("J110").Value = 1981.25
$ws.Range("L110").Value = 1981.25
$ws.Range("N110").Value = -6071.25
$ws.Range("H136").Value = 4010.2712
$ws.Range("I136").Value = 5068.4414
$ws.Range("K136").Value = 15205.3242
$ws.Range("M136").Value = -12655.3242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 11448578
$ws.Range("I19").Value = 20000
$ws.Range("J19").Value = 13353342
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 13353342
$ws.Range("M19").Value = -19827
$ws.Range("N19").Value = -13353688
$ws.Range("H26").Value = 7465.5
$ws.Range("I26").Value = 7465.5
$ws.Range("K26").Value = 7465.5
$ws.Range("M26").Value = -7173.5
$ws.Range("H99").Value = 1323.75
$ws.Range("I99").Value = 755.8333
$ws.Range("J99").Value = 1891.6666
$ws.Range("K99").Value = 755.8333
$ws.Range("L99").Value = 1891.6666
$ws.Range("M99").Value = 742.1667
$ws.Range("N99").Value = -4887.6666
$ws.Range("H107").Value = 875.05884
$ws.Range("I107").Value = 919.3
$ws.Range("J107").Value = 811.8570999999999
$ws.Range("K107").Value = 919.3
$ws.Range("L107").Value = 811.8570999999999
$ws.Range("M107").Value = 1000.7
$ws.Range("N107").Value = -4651.8571
$ws.Range("H134").Value = 1792.7609
$ws.Range("I134").Value = 1190.3611
$ws.Range("K134").Value = 3571.0833
$ws.Range("M134").Value = -1036.0833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3580.79
$ws.Range("I31").Value = 898.4545000000001
$ws.Range("J31").Value = 5424.896
$ws.Range("K31").Value = 898.4545000000001
$ws.Range("L31").Value = 5424.896
$ws.Range("M31").Value = -603.4545000000001
$ws.Range("N31").Value = -6014.896
$ws.Range("H34").Value = 3580.79
$ws.Range("I34").Value = 898.4545000000001
$ws.Range("J34").Value = 5424.896
$ws.Range("K34").Value = 898.4545000000001
$ws.Range("L34").Value = 5424.896
$ws.Range("M34").Value = -696.4545000000001
$ws.Range("N34").Value = -5828.896
$ws.Range("H58").Value = 1504.9474
$ws.Range("I58").Value = 1341.1666
$ws.Range("J58").Value = 1785.7142
$ws.Range("K58").Value = 1341.1666
$ws.Range("L58").Value = 1785.7142
$ws.Range("M58").Value = -1138.1666
$ws.Range("N58").Value = -2191.7142
$ws.Range("H132").Value = 2921.7144
$ws.Range("I132").Value = 2024
$ws.Range("J132").Value = 4118.6665
$ws.Range("K132").Value = 6072
$ws.Range("L132").Value = 12355.9995
$ws.Range("M132").Value = -3542
$ws.Range("N132").Value = -17415.9995
$ws.Range("H134").Value = 5075.154
$ws.Range("I134").Value = 5414.75
$ws.Range("K134").Value = 16244.25
$ws.Range("M134").Value = -13709.25
$ws.Range("H136").Value = 1504.9474
$ws.Range("I136").Value = 1341.1666
$ws.Range("J136").Value = 1785.7142
$ws.Range("K136").Value = 4023.4998
$ws.Range("L136").Value = 5357.142599999999
$ws.Range("M136").Value = -1473.4998
$ws.Range("N136").Value = -10457.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1029.34
$ws.Range("I68").Value = 757.6977000000001
$ws.Range("J68").Value = 1234.2632
$ws.Range("K68").Value = 2273.0931
$ws.Range("L68").Value = 3702.7896
$ws.Range("M68").Value = -1462.0931
$ws.Range("N68").Value = -5324.7896
$ws.Range("H71").Value = 1029.34
$ws.Range("I71").Value = 757.6977000000001
$ws.Range("J71").Value = 1234.2632
$ws.Range("K71").Value = 6819.2793
$ws.Range("L71").Value = 11108.3688
$ws.Range("M71").Value = -2763.2793
$ws.Range("N71").Value = -19220.3688
$ws.Range("H74").Value = 200
$ws.Range("I74").Value = 200
$ws.Range("K74").Value = 600
$ws.Range("M74").Value = 461
$ws.Range("H77").Value = 200
$ws.Range("I77").Value = 200
$ws.Range("K77").Value = 1800
$ws.Range("M77").Value = 3504
$ws.Range("H81").Value = 4403.6
$ws.Range("J81").Value = 5376.25
$ws.Range("L81").Value = 16128.75
$ws.Range("N81").Value = -18374.75
$ws.Range("H84").Value = 4403.6
$ws.Range("J84").Value = 5376.25
$ws.Range("L84").Value = 48386.25
$ws.Range("N84").Value = -59618.25
$ws.Range("H122").Value = 430.05264
$ws.Range("J122").Value = 714.8
$ws.Range("L122").Value = 6433.2
$ws.Range("N122").Value = -11333.2
$ws.Range("H132").Value = 1461.7587
$ws.Range("I132").Value = 1555.4117
$ws.Range("J132").Value = 1329.0834
$ws.Range("K132").Value = 13998.7053
$ws.Range("L132").Value = 11961.7506
$ws.Range("M132").Value = -11468.7053
$ws.Range("N132").Value = -17021.7506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1890
$ws.Range("I113").Value = 1880
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1880
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 290
$ws.Range("N113").Value = -6240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2478.762
$ws.Range("I61").Value = 2091.4119
$ws.Range("K61").Value = 2091.4119
$ws.Range("M61").Value = -1889.4119
$ws.Range("H113").Value = 2478.762
$ws.Range("I113").Value = 2091.4119
$ws.Range("K113").Value = 2091.4119
$ws.Range("M113").Value = 78.58809999999994
$ws.Range("H122").Value = 5052396.5
$ws.Range("I122").Value = 11112902
$ws.Range("J122").Value = 1975.8334
$ws.Range("K122").Value = 33338706
$ws.Range("L122").Value = 5927.5002
$ws.Range("M122").Value = -33336256
